$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data (shifts rows 1-2 down to 2-3)
$ws.Range("A1:C1").Insert()

# Populate the new header row
$ws.Range("A1").Value = "Supplier ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"

# Update selection to match the target state
$ws.Range("C1").Select()
